# Daily attendance processing - 2025-12-25 19:06:34
# Swap the first two comma-separated entries in the "Recorded By" column (G)
# for every data row on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = 7
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ", "

        if ($parts.Count -ge 2) {
            if ($parts.Count -gt 2) {
                $rest = $parts[2..($parts.Count - 1)]
            } else {
                $rest = @()
            }
            $swapped = @($parts[1], $parts[0]) + $rest
            $newValue = $swapped -join ", "
            $cell.Value2 = $newValue
        }
    }
}
